$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# Update the MAY (column H) budget values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2

# Move the active selection to H4 (matches the saved cursor position).
$ws.Range("H4").Select()
